# Adds a new "Quantidade Sugerida" (suggested quantity) column to the
# requisition header block, shifting the existing "Quantidade Requisitada",
# "Quantidade Autorizada" and "Quantidade Aprovada" columns one cell to the
# right (J->K, K->L, L->M) in both the label row (6) and the placeholder
# row (7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header labels (row 6) and placeholder tokens (row 7) one
# column to the right: L<-K, K<-J (do the rightmost column first so we
# do not clobber a value before it has been copied). Value2 is used
# because Value's getter is unreliable in this host.
$ws.Range("M6").Value2 = $ws.Range("L6").Value2
$ws.Range("L6").Value2 = $ws.Range("K6").Value2
$ws.Range("K6").Value2 = $ws.Range("J6").Value2
$ws.Range("J6").Value2 = " Quantidade Sugerida"

$ws.Range("M7").Value2 = $ws.Range("L7").Value2
$ws.Range("L7").Value2 = $ws.Range("K7").Value2
$ws.Range("K7").Value2 = $ws.Range("J7").Value2
$ws.Range("J7").Value2 = "{product.suggestedQuantity}"

# The shifted-into cells (M6/M7) were previously blank placeholders without
# the header border/centering; match the look of their new siblings by
# copying just the formatting over (values are left untouched).
$ws.Range("L6").Copy() | Out-Null
$ws.Range("M6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("L7").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active selection the same way the author left it.
$ws.Range("K9").Select()
